$d = $word.ActiveDocument

$oldText = "Chatting Socket IO"
$newText = "Chatlog"

# Word tracks the single most recent edit location with a hidden
# "_GoBack" bookmark. Right now it sits at the very end of the document
# (right after the last screenshot). Drop it here - it gets re-created
# below, collapsed at the *new* edit location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the typo without touching it yet, so the bookmark can be
# parked exactly where the fixed text is about to land.
$target = $d.Content
$found = $target.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if ($found) {
    $target.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $target) | Out-Null
}

# Now apply the actual fix: "Chatting Socket IO" -> "Chatlog".
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2) | Out-Null
